$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update absorption cross-section values (C4, C5) from 0.5 to 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

# Update the active selection to E3
$ws.Range("E3").Select()
